$d = $word.ActiveDocument

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# Locate the table. NOTE: materializing a table item straight off
# $d.Tables leaves $d.Paragraphs indexing stale/wrong afterwards in this
# host, so go through a full-document Range's Tables collection instead
# - that resolves the same table without disturbing $d.Paragraphs.
$fullRange = $d.Range(0, $d.Content.End)
$tbl = $fullRange.Tables.Item(1)
$tblStart = $tbl.Range.Start
$tblEnd = $tbl.Range.End

# --- Before the table ----------------------------------------------------
# The empty paragraph immediately before the <w:tbl> currently holds a
# single run with an empty <w:t/>: <w:p><w:r><w:t/></w:r></w:p>. Its
# paragraph mark is the one character position right before the table
# starts. Replacing that whole range (mark included) with a bare
# paragraph fragment drops the run but keeps the paragraph (and its
# mark) intact, so no paragraphs get merged/lost - result: <w:p/>.
$beforeTablePara = $d.Range($tblStart - 1, $tblStart)
[void]$beforeTablePara.InsertXML("<w:p xmlns:w='$wNs'/>")

# --- After the table -------------------------------------------------------
# The first empty paragraph right after the table is a bare <w:p/>; its
# start coincides with where the table's range ends. Give it a run with
# an empty <w:t/>, turning it into <w:p><w:r><w:t/></w:r></w:p>, by
# inserting that fragment at that (collapsed) point.
$afterTablePoint = $d.Range($tblEnd, $tblEnd)
[void]$afterTablePoint.InsertXML("<w:p xmlns:w='$wNs'><w:r><w:t></w:t></w:r></w:p>")

Write-Output "done"
